# Update "想去人数" (want-to-go count) column F for rows 2-13
# on both the "展览" and "全部类型" worksheets (they hold duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 417
    3  = 1412
    4  = 7158
    5  = 529
    6  = 297
    7  = 4935
    8  = 124
    9  = 1643
    10 = 53
    11 = 992
    12 = 258
    13 = 5447
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
